$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.62"
$ws.Range("E2").Value = "'1.06%"
$ws.Range("D3").Value = "'31.92"
$ws.Range("E3").Value = "'0.48%"
$ws.Range("D4").Value = "'5.027"
$ws.Range("E4").Value = "'-1.24%"
$ws.Range("D5").Value = "'0.07845"
$ws.Range("E5").Value = "'-3.77%"
$ws.Range("D6").Value = "'2.118"
$ws.Range("E6").Value = "'-16.03%"
$ws.Range("D7").Value = "'7.801"
$ws.Range("E7").Value = "'0.40%"
$ws.Range("D8").Value = "'3.783"
$ws.Range("E8").Value = "'-1.64%"
$ws.Range("D9").Value = "'0.9260"
$ws.Range("E9").Value = "'-0.46%"
$ws.Range("D10").Value = "'0.1742"
$ws.Range("E10").Value = "'-0.88%"
$ws.Range("D11").Value = "'0.07926"
$ws.Range("E11").Value = "'5.44%"
$ws.Range("D12").Value = "'0.08786"
$ws.Range("E12").Value = "'-2.07%"
$ws.Range("D13").Value = "'0.03127"
$ws.Range("E13").Value = "'4.35%"
$ws.Range("D14").Value = "'0.1002"
$ws.Range("E14").Value = "'0.18%"
$ws.Range("D15").Value = "'0.001508"
$ws.Range("E15").Value = "'0.16%"
$ws.Range("D16").Value = "'0.005956"
$ws.Range("E16").Value = "'3.39%"
$ws.Range("D17").Value = "'3.452"
$ws.Range("E17").Value = "'-3.61%"
$ws.Range("D18").Value = "'2.268"
$ws.Range("E18").Value = "'0.65%"
$ws.Range("D19").Value = "'0.3274"
$ws.Range("E20").Value = "'-2.74%"
$ws.Range("D21").Value = "'4.146"
$ws.Range("E21").Value = "'5.60%"
$ws.Range("E22").Value = "'5.49%"
$ws.Range("D23").Value = "'0.04601"
$ws.Range("E23").Value = "'0.06%"
$ws.Range("D24").Value = "'0.001236"
$ws.Range("E24").Value = "'-0.33%"
$ws.Range("D25").Value = "'0.004482"
$ws.Range("E25").Value = "'0.54%"
$ws.Range("D26").Value = "'0.0001248"
$ws.Range("E26").Value = "'4.11%"
$ws.Range("D39").Value = "'0.01737"
$ws.Range("E39").Value = "'-1.55%"
$ws.Range("D40").Value = "'0.04757"
$ws.Range("E40").Value = "'4.85%"
$ws.Range("D41").Value = "'0.007273"
$ws.Range("E41").Value = "'5.40%"
$ws.Range("D42").Value = "'0.1360"
$ws.Range("E42").Value = "'0.45%"
$ws.Range("D43").Value = "'0.002077"
$ws.Range("E43").Value = "'-5.90%"
$ws.Range("D44").Value = "'0.01078"
$ws.Range("E44").Value = "'10.21%"
$ws.Range("D45").Value = "'0.00006074"
$ws.Range("E45").Value = "'-7.35%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'0.10%"
$ws.Range("D47").Value = "'0.003399"
$ws.Range("E47").Value = "'-61.10%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.10%"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'0.10%"
